# Update "想去人数" (F column) and one "最低票价" (G21) value on both the
# "展览" sheet and the "全部类型" sheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value to apply on each of the two sheets.
$updates = @{
    "F2"  = 222
    "F4"  = 12967
    "F5"  = 1337
    "F6"  = 214
    "F7"  = 37
    "F8"  = 96
    "F10" = 226
    "F11" = 469
    "F12" = 3
    "F17" = 410
    "F18" = 5522
    "F20" = 52
    "F21" = 960
    "G21" = 58
    "F22" = 33
    "F23" = 130
    "F24" = 136
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
